$wb = $excel.ActiveWorkbook

# --- Data Entry sheet: change selection from E6 to A1:H4 ---
$wsData = $wb.Worksheets.Item("Data Entry")
$wsData.Range("A1:H4").Select()

# --- DateTime sheet: change selection from B15 to B33 ---
$wsDate = $wb.Worksheets.Item("DateTime")
$wsDate.Range("B33").Select()

# --- Add the new "Formulas" sheet after DateTime ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsFormulas = $wb.Worksheets.Add($null, $lastSheet)
$wsFormulas.Name = "Formulas"

# Bring over the Sales/Expenses/Profits table from Data Entry (A1:G4)
$wsData.Range("A1:G4").Copy()
$wsFormulas.Range("A1").PasteSpecial()

# Difference formula (written first so its label claims the next shared-string slot)
$wsFormulas.Range("B4").Formula = "=B2-B3"
$wsFormulas.Range("B5").Value = "Formula is =B2-B3 "

# Sum / Average column headers
$wsFormulas.Range("H1").Value = "Sum"
$wsFormulas.Range("I1").Value = "Average"

# Apply the new font (size 12, black) - once across the header row (A1:I1),
# then across the A1:H4 block (mirrors the author's formatting passes)
$wsFormulas.Range("A1:I1").Font.Size = 12
$wsFormulas.Range("A1:I1").Font.Color = 0
$wsFormulas.Range("A1:H4").Font.Size = 12
$wsFormulas.Range("A1:H4").Font.Color = 0

# Sum / Average formulas
$wsFormulas.Range("H2").Formula = "=SUM(B2:G2)"
$wsFormulas.Range("I2").Formula = "=AVERAGE(B2:G2)"
$wsFormulas.Range("H5").Value = "Sum formula is =SUM(B2:G2)"
$wsFormulas.Range("I5").Value = "Avg formula is =AVERAGE(B2:G2)"

# Widen H/I so the long labels fit (close match to the authored best-fit widths)
$wsFormulas.Columns("H:H").ColumnWidth = 24.998697916666668
$wsFormulas.Columns("I:I").ColumnWidth = 28.498697916666668

# Leave the cursor where the author left it
$wsFormulas.Range("I9").Select()
